# various fixes & improvement to live simulations + logging
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: raw poll numbers ---
$ws.Range("A2").Value = 26
$ws.Range("B2").Value = 35
$ws.Range("C2").Value = 13
$ws.Range("I2").Value = 14

# --- Row 10: second poll pair ---
$ws.Range("A10").Value = 45
$ws.Range("B10").Value = 49

# --- Row 24: live sim raw figures ---
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 34
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 16

# I24 previously held a formula (=1.89+0.6+4.68); it's now a plain literal value
$ws.Range("I24").Value = 13

# --- Row 25 ---
$ws.Range("C25").Value = 66.8
$ws.Range("D25").Value = 43
$ws.Range("I25").Value = 43

# --- Row 26 ---
$ws.Range("C26").Value = 89.1
$ws.Range("D26").Value = 25.5
$ws.Range("I26").Value = 52.8

# move the selection cursor to K8, matching the saved view state
$ws.Range("K8").Select()
